# Refresh cryptocurrency market snapshot (price / 1h volume / occasional rank-neighbour swaps)
# sourced from coinranking.com, as produced by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry maps a row to the columns that changed since the last run.
$updates = @(
    @{ Row=2; D='61.958.53'; E='  -1.14%  ' }
    @{ Row=3; D='3.405.46'; E='  -2.07%  ' }
    @{ Row=4; E='  +0.03%  ' }
    @{ Row=5; D='407.46'; E='  -1.33%  ' }
    @{ Row=6; D='133.48'; E='  +3.48%  ' }
    @{ Row=7; D='0.593'; E='  -0.89%  ' }
    @{ Row=8; E='  -0.04%  ' }
    @{ Row=9; D='0.682'; E='  -1.95%  ' }
    @{ Row=10; E='  -5.84%  ' }
    @{ Row=11; D='42.64'; E='  -2.29%  ' }
    @{ Row=13; D='8.40'; E='  -4.02%  ' }
    @{ Row=14; D='19.84'; E='  -2.20%  ' }
    @{ Row=15; D='3.439.50'; E='  -3.99%  ' }
    @{ Row=16; D='62.027.56'; E='  -0.83%  ' }
    @{ Row=17; E='  -3.28%  ' }
    @{ Row=18; D='11.00'; E='  -1.64%  ' }
    @{ Row=19; E='  -5.72%  ' }
    @{ Row=20; E='  -5.31%  ' }
    @{ Row=21; D='84.02'; E='  +2.03%  ' }
    @{ Row=22; D='313.45'; E='  +0.24%  ' }
    @{ Row=23; D='12.84'; E='  -2.91%  ' }
    @{ Row=24; E='  -0.91%  ' }
    @{ Row=25; E='  +9.76%  ' }
    @{ Row=26; D='29.55'; E='  -2.90%  ' }
    @{ Row=27; D='8.15'; E='  -0.27%  ' }
    @{ Row=28; D='2.80'; E='  +4.26%  ' }
    @{ Row=29; E='  -3.59%  ' }
    @{ Row=30; E='  -2.44%  ' }
    @{ Row=31; D='0.115'; E='  -4.05%  ' }
    @{ Row=32; D='42.76'; E='  -5.00%  ' }
    @{ Row=33; E='  -0.21%  ' }
    @{ Row=34; D='11.38'; E='  -6.48%  ' }
    @{ Row=35; D='0.0482'; E='  -2.82%  ' }
    @{ Row=36; D='51.69'; E='  -1.62%  ' }
    @{ Row=37; D='1.00'; E='  +0.37%  ' }
    @{ Row=38; E='  -4.73%  ' }
    @{ Row=39; E='  -3.38%  ' }
    @{ Row=40; E='  -0.76%  ' }
    @{ Row=41; B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.300'; E='  +3.53%  ' }
    @{ Row=42; D='137.24'; E='  -0.33%  ' }
    @{ Row=43; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.125'; E='  -0.81%  ' }
    @{ Row=44; D='4.03'; E='  +0.14%  ' }
    @{ Row=45; D='16.70'; E='  -7.01%  ' }
    @{ Row=46; D='2.22'; E='  -2.82%  ' }
    @{ Row=47; D='21.28'; E='  -5.56%  ' }
    @{ Row=48; D='2.120.04'; E='  -4.66%  ' }
    @{ Row=49; E='  -3.75%  ' }
    @{ Row=50; E='  +2.57%  ' }
    @{ Row=51; E='  +16.43%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in 'B','C','D','E') {
        if (-not $u.ContainsKey($col)) { continue }
        $cell = $ws.Range("$col$row")
        if ($col -eq 'D') {
            # Price column holds text (dotted-thousands, trailing zeros) -- keep it text
            # instead of letting Excel coerce it to a trimmed number.
            $cell.NumberFormat = "@"
            $cell.Value = $u[$col]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u[$col]
        }
    }
}
